$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.00%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.520"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.39%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08473"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.66%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.069"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.40%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9885"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.542"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.43%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1176"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.86%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1923"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.99%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "9.501"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.90%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09836"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.25%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04706"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.71%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1060"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.24%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001284"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.15%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005902"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.87%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.21%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.442"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.76%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3330"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.84%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1386"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.18%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2551"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.13%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04160"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.83%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001304"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.23%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004598"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "5.73%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001304"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.63%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002989"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.23%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02709"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "5.17%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05754"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.32%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007803"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.92%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1433"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.19%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007445"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.52%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002159"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.27%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008071"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3555"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007071"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.55%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.26%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.21%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003448"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-1.57%"
$ws.Range("B50").Value = "CoinbaseStockToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.003542"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.18%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002107"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.26%"
